$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.21466851234436
$ws.Range("B1").Value = 2.599881649017334
$ws.Range("C1").Value = 4.383058547973633
$ws.Range("D1").Value = 2.129381418228149
$ws.Range("E1").Value = 1.166103959083557
